$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name cell (A column) where row ranking changed ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 18:05"
$ws.Range("A55").Value = "Argelia"
$ws.Range("A56").Value = "Barein"
$ws.Range("A62").Value = "Moldavia"
$ws.Range("A63").Value = "Nigeria"
$ws.Range("A69").Value = "Irak"
$ws.Range("A70").Value = "Hungria"
$ws.Range("A196").Value = "Nueva Caledonia"
$ws.Range("A197").Value = "Belice"
$ws.Range("A209").Value = "Seychelles"
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Montserrat"
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# --- Update statistic cells (B-H columns) with refreshed case counts ---
$ws.Range("B4").Value = 1555768
$ws.Range("C4").Value = 5474
$ws.Range("D4").Value = 359175
$ws.Range("E4").Value = 1104308
$ws.Range("G4").Value = 304
$ws.Range("H4").Value = 92285
$ws.Range("B7").Value = 261567
$ws.Range("C7").Value = 6199
$ws.Range("E7").Value = 143733
$ws.Range("G7").Value = 522
$ws.Range("H7").Value = 17375
$ws.Range("B14").Value = 103886
$ws.Range("C14").Value = 3558
$ws.Range("D14").Value = 40856
$ws.Range("E14").Value = 59818
$ws.Range("G14").Value = 56
$ws.Range("H14").Value = 3212
$ws.Range("B17").Value = 78499
$ws.Range("C17").Value = 427
$ws.Range("D17").Value = 39488
$ws.Range("E17").Value = 33154
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 5857
$ws.Range("B21").Value = 49579
$ws.Range("C21").Value = 3520
$ws.Range("D21").Value = 21507
$ws.Range("E21").Value = 27563
$ws.Range("G21").Value = 31
$ws.Range("H21").Value = 509
$ws.Range("E28").Value = 1127
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 1891
$ws.Range("B34").Value = 19268
$ws.Range("C34").Value = 383
$ws.Range("E34").Value = 10417
$ws.Range("G34").Value = 12
$ws.Range("H34").Value = 948
$ws.Range("E37").Value = 5888
$ws.Range("G37").Value = 17
$ws.Range("H37").Value = 1137
$ws.Range("B55").Value = 7377
$ws.Range("C55").Value = 176
$ws.Range("D55").Value = 3625
$ws.Range("E55").Value = 3191
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 561
$ws.Range("B56").Value = 7374
$ws.Range("C56").Value = 190
$ws.Range("D56").Value = 2952
$ws.Range("E56").Value = 4410
$ws.Range("H56").Value = 12
$ws.Range("D60").Value = 3598
$ws.Range("E60").Value = 3118
$ws.Range("B62").Value = 6340
$ws.Range("C62").Value = 202
$ws.Range("D62").Value = 2508
$ws.Range("E62").Value = 3611
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 221
$ws.Range("B63").Value = 6175
$ws.Range("D63").Value = 1644
$ws.Range("E63").Value = 4340
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 191
$ws.Range("B69").Value = 3611
$ws.Range("C69").Value = 57
$ws.Range("D69").Value = 2366
$ws.Range("E69").Value = 1114
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 131
$ws.Range("B70").Value = 3556
$ws.Range("C70").Value = 21
$ws.Range("D70").Value = 1412
$ws.Range("E70").Value = 1677
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 467
$ws.Range("B74").Value = 2840
$ws.Range("C74").Value = 4
$ws.Range("E74").Value = 1301
$ws.Range("B110").Value = 918
$ws.Range("C110").Value = 1
$ws.Range("E110").Value = 386
$ws.Range("B133").Value = 420
$ws.Range("C133").Value = 8
$ws.Range("D133").Value = 132
$ws.Range("E133").Value = 273
$ws.Range("B135").Value = 391
$ws.Range("C135").Value = 3
$ws.Range("E135").Value = 43
$ws.Range("D196").Value = 18
$ws.Range("H196").Value = 0
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
